$wb = $excel.ActiveWorkbook

# --- Report regeneration: bump localization status from "Ready for handoff"
# --- to "In Translation" everywhere it is surfaced, then re-fit the status
# --- columns now that the label text is shorter.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Column width is stored (OOXML) as characters = ColumnWidth + ~0.8333 padding,
# so feed back the padding-compensated figure to land on the report's target
# rendered width for the (now narrower) status columns.
$targetColumnWidth = 12.576851254417766

$wsOverview.Range("E:E").ColumnWidth = $targetColumnWidth
$wsOverview.Range("F:F").ColumnWidth = $targetColumnWidth
$wsZhCn.Range("C:C").ColumnWidth = $targetColumnWidth
$wsDeDe.Range("C:C").ColumnWidth = $targetColumnWidth
